$d = $word.ActiveDocument

# 1) Update the course term in the title paragraph: Fall 2019 -> Spring 2021.
#    Locate the paragraph by its current text and rewrite the whole range so
#    that the run keeps its original formatting (bold) and xml:space handling.
foreach ($p in $d.Paragraphs) {
  if ($p.Range.Text -like "*Educational Technology 0858-612, Fall 2019*") {
    $p.Range.Text = "Educational Technology 0858-612, Spring 2021"
    break
  }
}

# 2) Insert a new "Spring 2021" weekly schedule table (week / date / format)
#    immediately after the existing weekly-topics table (Table 1), i.e.
#    right before the "Assignments & grading" heading.
$t1 = $d.Tables.Item(1)
$insertPos = $t1.Range.End
$insertRange = $d.Range($insertPos, $insertPos)

$newTableXml = @'
<w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:tblPr><w:tblStyle w:val="Table"/><w:tblW w:type="pct" w:w="0.0"/><w:tblLook w:firstRow="1"/></w:tblPr><w:tblGrid/><w:tr><w:trPr><w:cnfStyle w:firstRow="1"/></w:trPr><w:tc><w:tcPr><w:tcBorders><w:bottom w:val="single"/></w:tcBorders><w:vAlign w:val="bottom"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">week</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcBorders><w:bottom w:val="single"/></w:tcBorders><w:vAlign w:val="bottom"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">date</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcBorders><w:bottom w:val="single"/></w:tcBorders><w:vAlign w:val="bottom"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">format</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">1</w:t></w:r></w:p></w:tc><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">Jan 27 Going mobile</w:t></w:r></w:p></w:tc><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">flex</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">2</w:t></w:r></w:p></w:tc><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">Feb 03 Mobile first</w:t></w:r></w:p></w:tc><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">online</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">3</w:t></w:r></w:p></w:tc><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">Feb 10 Tech reports</w:t></w:r></w:p></w:tc><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">flex</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">4</w:t></w:r></w:p></w:tc><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">Feb 17 Mobile computing and society</w:t></w:r></w:p></w:tc><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">online</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">5</w:t></w:r></w:p></w:tc><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">Feb 24 App Inventor hackathon</w:t></w:r></w:p></w:tc><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">flex</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">6</w:t></w:r></w:p></w:tc><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">Mar 03 Mobile cognition</w:t></w:r></w:p></w:tc><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">online</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">-</w:t></w:r></w:p></w:tc><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">Mar 10</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">no class (mini break)</w:t></w:r></w:p></w:tc><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">-</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">7</w:t></w:r></w:p></w:tc><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">Mar 17 Augmented reality</w:t></w:r></w:p></w:tc><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">flex</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">8</w:t></w:r></w:p></w:tc><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">Mar 24 1:1 Computing</w:t></w:r></w:p></w:tc><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">online</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">9</w:t></w:r></w:p></w:tc><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">Mar 31 Reading screens</w:t></w:r></w:p></w:tc><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">flex</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">10</w:t></w:r></w:p></w:tc><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">Apr 07 Mobile games for learning</w:t></w:r></w:p></w:tc><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">online</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">11</w:t></w:r></w:p></w:tc><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">Apr 14 Workshop: ARIS</w:t></w:r></w:p></w:tc><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">flex</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">12</w:t></w:r></w:p></w:tc><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">Apr 21 Internet of Things (IoT)</w:t></w:r></w:p></w:tc><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">online</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">13</w:t></w:r></w:p></w:tc><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">Apr 28 Mobile design (UX)</w:t></w:r></w:p></w:tc><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">flex</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">14</w:t></w:r></w:p></w:tc><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">May 05 Final project workshop</w:t></w:r></w:p></w:tc><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">online</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">-</w:t></w:r></w:p></w:tc><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">May 12</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">no class (makeup day)</w:t></w:r></w:p></w:tc><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">-</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">15</w:t></w:r></w:p></w:tc><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">May 19 Final project presentation</w:t></w:r></w:p></w:tc><w:tc><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">flex (final presentations)</w:t></w:r></w:p></w:tc></w:tr></w:tbl>
'@

$insertRange.InsertXML($newTableXml) | Out-Null
